$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1961.0769
$ws.Range("I40").Value = 1846.3077
$ws.Range("J40").Value = 2075.8462
$ws.Range("K40").Value = 1846.3077
$ws.Range("L40").Value = 2075.8462
$ws.Range("M40").Value = -1671.3077
$ws.Range("N40").Value = -2425.8462
$ws.Range("H51").Value = 3243.2856
$ws.Range("J51").Value = 3300.5
$ws.Range("L51").Value = 3300.5
$ws.Range("N51").Value = -4268.5
$ws.Range("H58").Value = 934.93335
$ws.Range("J58").Value = 3005.6667
$ws.Range("L58").Value = 9017.000100000001
$ws.Range("N58").Value = -9317.000100000001
$ws.Range("H86").Value = 1814.8
$ws.Range("I86").Value = 2251.875
$ws.Range("J86").Value = 1037.7778
$ws.Range("K86").Value = 2251.875
$ws.Range("L86").Value = 1037.7778
$ws.Range("M86").Value = -1128.875
$ws.Range("N86").Value = -3283.7778
$ws.Range("H89").Value = 1814.8
$ws.Range("I89").Value = 2251.875
$ws.Range("J89").Value = 1037.7778
$ws.Range("K89").Value = 11259.375
$ws.Range("L89").Value = 5188.889
$ws.Range("M89").Value = -5643.375
$ws.Range("N89").Value = -16420.889
$ws.Range("H137").Value = 1463.8108
$ws.Range("I137").Value = 1135.3
$ws.Range("J137").Value = 1850.2941
$ws.Range("K137").Value = 3405.9
$ws.Range("L137").Value = 5550.8823
$ws.Range("M137").Value = -855.8999999999996
$ws.Range("N137").Value = -10650.8823

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2472.18
$ws.Range("I32").Value = 1992.4713
$ws.Range("J32").Value = 5682.5386
$ws.Range("K32").Value = 1992.4713
$ws.Range("L32").Value = 5682.5386
$ws.Range("M32").Value = -1705.4713
$ws.Range("N32").Value = -6256.5386
$ws.Range("H61").Value = 76924440
$ws.Range("I61").Value = 90910250
$ws.Range("J61").Value = 2507
$ws.Range("K61").Value = 90910250
$ws.Range("L61").Value = 2507
$ws.Range("M61").Value = -90910038
$ws.Range("N61").Value = -2931
$ws.Range("H107").Value = 30001
$ws.Range("J107").Value = 30001
$ws.Range("L107").Value = 30001
$ws.Range("N107").Value = -37681
$ws.Range("H122").Value = 2938.5
$ws.Range("I122").Value = 2737.5
$ws.Range("K122").Value = 8212.5
$ws.Range("M122").Value = -5762.5
$ws.Range("H136").Value = 76924440
$ws.Range("I136").Value = 90910250
$ws.Range("J136").Value = 2507
$ws.Range("K136").Value = 272730750
$ws.Range("L136").Value = 7521
$ws.Range("M136").Value = -272728200
$ws.Range("N136").Value = -12621

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2581.8
$ws.Range("J86").Value = 2843.25
$ws.Range("L86").Value = 2843.25
$ws.Range("N86").Value = -5089.25
$ws.Range("H89").Value = 2581.8
$ws.Range("J89").Value = 2843.25
$ws.Range("L89").Value = 14216.25
$ws.Range("N89").Value = -25448.25
$ws.Range("H105").Value = 1452.591
$ws.Range("I105").Value = 1459.7368
$ws.Range("J105").Value = 1407.3334
$ws.Range("K105").Value = 1459.7368
$ws.Range("L105").Value = 1407.3334
$ws.Range("M105").Value = 287.2632000000001
$ws.Range("N105").Value = -4901.3334
$ws.Range("H107").Value = 1493.2307
$ws.Range("I107").Value = 1088.7778
$ws.Range("J107").Value = 2403.25
$ws.Range("K107").Value = 1088.7778
$ws.Range("L107").Value = 2403.25
$ws.Range("M107").Value = 831.2221999999999
$ws.Range("N107").Value = -6243.25
$ws.Range("H134").Value = 1581.8182
$ws.Range("I134").Value = 1300.25
$ws.Range("J134").Value = 2332.6667
$ws.Range("K134").Value = 3900.75
$ws.Range("L134").Value = 6998.000100000001
$ws.Range("M134").Value = -1365.75
$ws.Range("N134").Value = -12068.0001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1060.6666
$ws.Range("I107").Value = 592.5714
$ws.Range("K107").Value = 592.5714
$ws.Range("M107").Value = 1327.4286
$ws.Range("H134").Value = 35717212
$ws.Range("I134").Value = 3881.2856
$ws.Range("K134").Value = 11643.8568
$ws.Range("M134").Value = -9108.856800000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H101").Value = 5005.625
$ws.Range("J101").Value = 5005.625
$ws.Range("L101").Value = 15016.875
$ws.Range("N101").Value = -19884.875
$ws.Range("H131").Value = 15152403
$ws.Range("J131").Value = 1039.0625
$ws.Range("L131").Value = 3117.1875
$ws.Range("N131").Value = -13197.1875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 965.86365
$ws.Range("I102").Value = 915.46155
$ws.Range("J102").Value = 1038.6666
$ws.Range("K102").Value = 915.46155
$ws.Range("L102").Value = 1038.6666
$ws.Range("M102").Value = 706.53845
$ws.Range("N102").Value = -4282.6666
$ws.Range("H107").Value = 587.06665
$ws.Range("I107").Value = 542.8570999999999
$ws.Range("K107").Value = 542.8570999999999
$ws.Range("M107").Value = 1377.1429
$ws.Range("H113").Value = 1530.8182
$ws.Range("I113").Value = 1369.875
$ws.Range("K113").Value = 1369.875
$ws.Range("M113").Value = 800.125
$ws.Range("H132").Value = 3531.5417
$ws.Range("I132").Value = 3430.8572
$ws.Range("J132").Value = 3672.5
$ws.Range("K132").Value = 10292.5716
$ws.Range("L132").Value = 11017.5
$ws.Range("M132").Value = -7762.571599999999
$ws.Range("N132").Value = -16077.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1200
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 1200
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 1200
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -1790
$ws.Range("H27").Value = 1200
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 1200
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 1200
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -1414
$ws.Range("H46").Value = 4066.5833
$ws.Range("I46").Value = 666.3333
$ws.Range("J46").Value = 5200
$ws.Range("K46").Value = 666.3333
$ws.Range("L46").Value = 5200
$ws.Range("M46").Value = -478.3333
$ws.Range("N46").Value = -5576
$ws.Range("H55").Value = 223.3158
$ws.Range("I55").Value = 102.52941
$ws.Range("J55").Value = 1250
$ws.Range("K55").Value = 102.52941
$ws.Range("L55").Value = 1250
$ws.Range("M55").Value = 70.47059
$ws.Range("N55").Value = -1596
$ws.Range("H132").Value = 2966.6191
$ws.Range("I132").Value = 2638.6155
$ws.Range("J132").Value = 3499.625
$ws.Range("K132").Value = 7915.8465
$ws.Range("L132").Value = 10498.875
$ws.Range("M132").Value = -5385.8465
$ws.Range("N132").Value = -15558.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1954.88
$ws.Range("I132").Value = 1741.2222
$ws.Range("J132").Value = 2504.2856
$ws.Range("K132").Value = 5223.6666
$ws.Range("L132").Value = 7512.8568
$ws.Range("M132").Value = -2693.6666
$ws.Range("N132").Value = -12572.8568
$ws.Range("H136").Value = 1265.0834
$ws.Range("I136").Value = 1119.125
$ws.Range("K136").Value = 3357.375
$ws.Range("M136").Value = -807.375
